# "Expense List Section" — update the existing income rows with new
# figures and append the previous rows further down the sheet, plus two
# new "Freelance" income rows, extending the table from A1:C5 to A1:C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Salary -> Interest From Saving Account, new amount/date ---
$ws.Range("A2").Value = "Interest From Saving Account"
$ws.Range("B2").Value = 4500
$ws.Range("C2").Value = 46024.229537037034

# --- Row 3: same source (Salary), new amount/date ---
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 46024.229537037034

# --- Row 4: same source (Salary), new amount/date ---
$ws.Range("B4").Value = 2400
$ws.Range("C4").Value = 46023.229537037034

# --- Row 5: same source (Salary), new amount/date ---
$ws.Range("B5").Value = 10000
$ws.Range("C5").Value = 46023.229537037034

# Carry the existing date formatting (style index used by C2:C5) down
# through the new rows C6:C10 before filling in their values, so the new
# date cells pick up the same number format as the rest of column C.
$ws.Range("C2").Copy()
$ws.Range("C6:C10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 6: new Freelance income ---
$ws.Range("A6").Value = "Freelance"
$ws.Range("B6").Value = 3000
$ws.Range("C6").Value = 46014.229537037034

# --- Row 7: previous row-2 data (Salary 5600) moved down ---
$ws.Range("A7").Value = "Salary"
$ws.Range("B7").Value = 5600
$ws.Range("C7").Value = 45728.229537037034

# --- Row 8: previous row-3 data (Salary 4600) moved down ---
$ws.Range("A8").Value = "Salary"
$ws.Range("B8").Value = 4600
$ws.Range("C8").Value = 45728.229537037034

# --- Row 9: previous row-4 data (Salary 4000) moved down ---
$ws.Range("A9").Value = "Salary"
$ws.Range("B9").Value = 4000
$ws.Range("C9").Value = 45728.229537037034

# --- Row 10: new Freelance income (replicates old row-5 date) ---
$ws.Range("A10").Value = "Freelance"
$ws.Range("B10").Value = 3000
$ws.Range("C10").Value = 45728.229537037034
